$d = $word.ActiveDocument

# Make sure inserted content isn't wrapped as tracked-change markup.
$d.TrackRevisions = $false

# Locate the unique "CANDY" entry in the red/bold keyword list
# ("..., ATTACK(...), CANDY, CORRELAT(...), ...").
$rng = $d.Content
$found = $rng.Find.Execute("CANDY", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'CANDY' in the document."
}

$candyStart = $rng.Start
$candyLen = $rng.End - $rng.Start

# ---------------------------------------------------------------------
# Insert "BANSHEE, " immediately before "CANDY".
# ---------------------------------------------------------------------
# Range.InsertBefore always lands right at the (fixed) start of the anchor
# range, so insert the later-appearing piece of text first; the text is
# typed adjacent to the plain (non-colored) run that already precedes
# CANDY, so it naturally inherits plain formatting.
$rng.InsertBefore(", ")
$rng.InsertBefore("BANSHEE")

# "BANSHEE" now occupies the 7 characters right before the original CANDY
# position; color/weight it like the other red keywords.
$bansheeRange = $d.Range($candyStart, $candyStart + 7)
$bansheeRange.Font.Bold = $true
$bansheeRange.Font.Color = 255

# The ", " separator between BANSHEE and CANDY stays plain bold (no
# explicit color - it already inherited the automatic color).
$sepBeforeCandy = $d.Range($candyStart + 7, $candyStart + 9)
$sepBeforeCandy.Font.Bold = $true

# Absolute position right after "CANDY" now that 9 characters were added
# in front of it.
$candyEndNow = $candyStart + 9 + $candyLen

# ---------------------------------------------------------------------
# Insert ", CAROSEL" immediately after "CANDY" (before "CORRELAT").
# ---------------------------------------------------------------------
# The original ", " that followed CANDY (two single-character runs, "," and
# " ") sits right here, immediately before "CORRELAT".
$oldSepEnd = $candyEndNow + 2

# Insert the new text adjacent to that plain old separator (i.e. right
# before CORRELAT) rather than directly after CANDY, so it inherits the
# plain (non-colored) formatting instead of CANDY's red.
$insPoint = $d.Range($oldSepEnd, $oldSepEnd)
$insPoint.InsertBefore(", ")
$insPoint.InsertBefore("CAROSEL")

# "CAROSEL": red bold, like the other keywords.
$caroselRange = $d.Range($oldSepEnd, $oldSepEnd + 7)
$caroselRange.Font.Bold = $true
$caroselRange.Font.Color = 255
